$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 5 holds the datetimes for the 33a4f6a5... handoff/handback pair
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-25 06:15:03"
$wsZhCn.Range("G5").Value = "2016-02-25 06:15:48"

# de-de sheet: row 5 holds the datetimes for the 33a4f6a5... handoff/handback pair
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-25 06:15:16"
$wsDeDe.Range("G5").Value = "2016-02-25 06:16:09"
